# Update the "Morning"/"Evening" measurement data table (columns B/C) for
# rows 76-84, adding the newly recorded values, and refresh the sheet's
# view/selection state to match where the user was working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (B value, C value); $null means "leave blank / do not set"
$newData = @{
    76 = @($null, 116.4)
    77 = @(115.5, 116.1)
    78 = @(115.2, 117)
    79 = @(115,   115.2)
    80 = @(114.6, 116.1)
    81 = @(115.1, 116)
    82 = @(115.2, 115.7)
    83 = @(114.6, 115.5)
    84 = @(114.4, $null)
}

foreach ($row in 76..84) {
    $pair = $newData[$row]
    $bVal = $pair[0]
    $cVal = $pair[1]

    if ($null -ne $bVal) {
        $ws.Cells.Item($row, 2).Value = $bVal
    }
    if ($null -ne $cVal) {
        $ws.Cells.Item($row, 3).Value = $cVal
    }
}

# Update the view: scroll position and active selection moved.
$ws.Activate() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollRow = 65
$win.ScrollColumn = 1
$ws.Range("B84").Select() | Out-Null

Write-Host "data updated"
